$d = $word.ActiveDocument

# Step 1: merge "Versi" + "on" runs into a single "Version" run.
# A temporary distinct value forces the engine to actually rewrite the
# run structure (setting identical text is treated as a no-op).
$r1 = $d.Range(0, 7)
$r1.Text = "VersionX"
$r1b = $d.Range(0, 8)
$r1b.Text = "Version"

# Step 2: change the " 2" run into " 1."
$r2 = $d.Range(7, 9)
$r2.Text = " 1."

# Step 3: delete the trailing "." run (now at offset 10-11)
$r3 = $d.Range(10, 11)
$r3.Delete()
